# Daily attendance processing - 2025-10-29 01:22:19
# Move "System" to the front of the "Recorded By" (column G) value list
# for every row on the active sheet, leaving everything else unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "

    # Case-sensitive check for an exact "System" token (comparisons with
    # -eq/-ne are case-insensitive in this engine, so use .Equals()).
    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
        $newText = $newParts -join ", "
        $cell.Value = $newText
    }
}
